$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values for row 2
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

# Move the active selection to D2 (was D5)
$ws.Range("D2").Select()
